# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# Re-groups the "Estado de Cuenta" detail rows (B16:G36) by worker instead
# of by period, and updates the "Valor Mora" (F column) for period 2412 so
# that each worker keeps one consistent value across periods 2406-2411 and
# a distinct (lower) value for period 2412.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Final data block for rows 16-36: TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico
$rows = @(
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2412", 22000, 1500000),
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2411", 60000, 1500000),
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2410", 60000, 1500000),
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2409", 60000, 1500000),
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2408", 60000, 1500000),
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2407", 60000, 1500000),
    @("CC", "1047419710", "KAREN PATRICIA TERAN GALLARDO", "2406", 60000, 1500000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2412", 19067, 1300000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2411", 52000, 1300000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2410", 52000, 1300000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2409", 52000, 1300000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2408", 52000, 1300000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2407", 52000, 1300000),
    @("CC", "73214409",   "JORGE ARMANDO MERCADO PATERNINA", "2406", 52000, 1300000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2412", 22000, 1500000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2411", 60000, 1500000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2410", 60000, 1500000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2409", 60000, 1500000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2408", 60000, 1500000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2407", 60000, 1500000),
    @("CC", "79598823",   "EDUARD TAMAYO RODRIGUEZ", "2406", 60000, 1500000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $data[1]   # C - N Doc Trabajador (text)
    $ws.Cells.Item($r, 4).Value = $data[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[3]   # E - Periodo Mora (text)
    $ws.Cells.Item($r, 6).Value = $data[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[5]   # G - Salario Basico
}
